# Add check-chord-type logic: if a chord's "chord type" (column G) is
# minor or diminished (e.g. "minor", "minor-seventh", "diminished", ...),
# bump its "chord encoded" value (column L) by 12 to flag it as such.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-IsMinorOrDiminished([string]$chordType) {
    if ([string]::IsNullOrEmpty($chordType)) { return $false }
    $t = $chordType.ToLower()
    return ($t -like "*minor*") -or ($t -like "*diminished*")
}

$lastRow = $ws.Cells.Item(1, 7).End(4).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $chordType = $ws.Cells.Item($row, 7).Text
    if (Test-IsMinorOrDiminished $chordType) {
        $encodedText = $ws.Cells.Item($row, 12).Text
        $encoded = [int]$encodedText
        $ws.Cells.Item($row, 12).Value = $encoded + 12
    }
}
